$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29 of the state-transition matrix was previously blank (style-only
# placeholder cells). The final state fills it in like the other symbol
# rows: column A gets the new "$" symbol label, and columns B..R get the
# "F" filler value used throughout the rest of the table.

# Column A: new symbol label, formatted like the other label cells (e.g. A4).
$ws.Range("A29").Value = "$"
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A29").PasteSpecial(-4122) | Out-Null

# Columns B..R: filler value "F", formatted like the other body cells (e.g. D4).
$ws.Range("B29:R29").Value = "F"
$ws.Range("D4").Copy() | Out-Null
$ws.Range("B29:R29").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0
